$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = $origStyle
}

Set-TextValue "D2" "51.491.05"
Set-TextValue "E2" "  +0.67%  "
Set-TextValue "D3" "2.973.09"
Set-TextValue "E3" "  +1.98%  "
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "378.37"
Set-TextValue "E5" "  +1.81%  "
Set-TextValue "D6" "104.34"
Set-TextValue "E6" "  +0.22%  "
Set-TextValue "D7" "0.541"
Set-TextValue "E7" "  +0.17%  "
Set-TextValue "D8" "1.00"
Set-TextValue "E8" "  +0.03%  "
Set-TextValue "D9" "0.592"
Set-TextValue "E9" "  +0.91%  "
Set-TextValue "D10" "37.20"
Set-TextValue "E10" "  +1.63%  "
Set-TextValue "E11" "  +0.19%  "
Set-TextValue "D12" "0.0843"
Set-TextValue "E12" "  +0.97%  "
Set-TextValue "D13" "3.440.91"
Set-TextValue "E13" "  +2.19%  "
Set-TextValue "D14" "18.42"
Set-TextValue "E14" "  +0.26%  "
Set-TextValue "D15" "7.57"
Set-TextValue "E15" "  +2.46%  "
Set-TextValue "D16" "2.975.11"
Set-TextValue "E16" "  +2.12%  "
Set-TextValue "D17" "0.966"
Set-TextValue "E17" "  +3.36%  "
Set-TextValue "D18" "51.443.77"
Set-TextValue "E18" "  +0.85%  "
Set-TextValue "E19" "  +2.44%  "
Set-TextValue "D20" "7.40"
Set-TextValue "E20" "  +2.65%  "
Set-TextValue "D21" "12.91"
Set-TextValue "E21" "  +0.15%  "
Set-TextValue "D22" "0.0₃0961"
Set-TextValue "E22" "  +1.87%  "
Set-TextValue "D23" "69.43"
Set-TextValue "E23" "  +1.74%  "
Set-TextValue "D24" "261.75"
Set-TextValue "E24" "  +0.90%  "
Set-TextValue "E25" "  +4.74%  "
Set-TextValue "D26" "8.23"
Set-TextValue "E26" "  +18.45%  "
Set-TextValue "D27" "7.54"
Set-TextValue "E27" "  +22.35%  "
Set-TextValue "E28" "  -0.54%  "
Set-TextValue "E29" "  +0.07%  "
Set-TextValue "E30" "  +9.07%  "
Set-TextValue "D31" "25.89"
Set-TextValue "E31" "  +0.71%  "
Set-TextValue "D32" "9.90"
Set-TextValue "E32" "  +0.12%  "
Set-TextValue "D33" "35.03"
Set-TextValue "E33" "  +1.00%  "
Set-TextValue "E34" "  -2.10%  "
Set-TextValue "D35" "50.99"
Set-TextValue "E35" "  +0.42%  "
Set-TextValue "E36" "  +5.83%  "
Set-TextValue "E37" "  +0.15%  "
Set-TextValue "D38" "3.04"
Set-TextValue "E38" "  +0.37%  "
Set-TextValue "D39" "17.16"
Set-TextValue "E39" "  +0.48%  "
Set-TextValue "D40" "2.58"
Set-TextValue "E40" "  -2.14%  "
Set-TextValue "D41" "1.85"
Set-TextValue "E41" "  +0.59%  "
Set-TextValue "D42" "0.116"
Set-TextValue "E42" "  +2.36%  "
Set-TextValue "D43" "125.01"
Set-TextValue "E43" "  +4.88%  "
Set-TextValue "D44" "21.68"
Set-TextValue "E44" "  -2.28%  "
Set-TextValue "D45" "0.290"
Set-TextValue "E45" "  +18.17%  "
Set-TextValue "E46" "  -1.64%  "
Set-TextValue "E47" "  +2.66%  "
Set-TextValue "D48" "2.035.16"
Set-TextValue "E48" "  +0.83%  "
Set-TextValue "D49" "3.22"
Set-TextValue "E49" "  +1.37%  "
Set-TextValue "D50" "0.0341"
Set-TextValue "E50" "  +10.10%  "
Set-TextValue "D51" "58.05"
Set-TextValue "E51" "  +2.31%  "
